# Adds a new "2022-Q3" quarterly sheet (9 funds) right after "总计",
# and inserts a matching summary row at the top of the "总计" sheet's
# data, shifting the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new data row right
#    below the header for the 2022-Q3 totals, pushing every other
#    quarter down by one row.
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()

# Match the look of the other index cells in column A (bold/centered
# style used by A3:A7).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 9
$summary.Range("D2").Value = 0.56

# ---------------------------------------------------------------
# 2. Insert the new "2022-Q3" sheet right after "总计" and before
#    the existing "2022-Q2" sheet.
# ---------------------------------------------------------------
$anchor = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $anchor)
$q3.Name = "2022-Q3"

# Header row (B1:H1), styled like the other quarterly sheets.
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
$headerCols = @(2,3,4,5,6,7,8)
for ($i = 0; $i -lt $headers.Length; $i++) {
  $cell = $q3.Cells.Item(1, $headerCols[$i])
  $cell.Value = $headers[$i]
}
$q3.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats, re-apply after value set
for ($i = 0; $i -lt $headers.Length; $i++) {
  $cell = $q3.Cells.Item(1, $headerCols[$i])
  $cell.Value = $headers[$i]
}

# Copy the header style (s="2", same as the other sheets) from the
# "总计" sheet's own header cells for a faithful match.
$summary.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
for ($i = 0; $i -lt $headers.Length; $i++) {
  $cell = $q3.Cells.Item(1, $headerCols[$i])
  $cell.Value = $headers[$i]
}

# Fund rows (row 2 .. row 10).
$funds = @(
  @("005314","万家中证1000指数增强C","14.28","94.11","1.03","0.1471",4),
  @("005313","万家中证1000指数增强A","13.25","94.11","1.03","0.1365",4),
  @("013641","博道成长智航股票A","10.27","90.29","0.84","0.0863",5),
  @("013642","博道成长智航股票C","7.24","90.29","0.84","0.0608",5),
  @("015784","中信建投中证1000指数增强A","8.10","92.20","0.71","0.0575",2),
  @("015785","中信建投中证1000指数增强C","3.32","92.20","0.71","0.0236",2),
  @("004194","招商中证1000指数增强A","1.56","92.06","1.18","0.0184",3),
  @("562900","易方达中证现代农业主题ETF","0.64","97.55","2.15","0.0138",9),
  @("004195","招商中证1000指数增强C","1.09","92.06","1.18","0.0129",3)
)

# Column A index cells share the bold/centered style (s="2") used in
# the other quarterly sheets.
$q3.Range("B1").Copy()
for ($i = 0; $i -lt $funds.Length; $i++) {
  $r = $i + 2
  $q3.Cells.Item($r, 1).PasteSpecial(-4122)  # xlPasteFormats
}

for ($i = 0; $i -lt $funds.Length; $i++) {
  $r = $i + 2
  $fund = $funds[$i]

  $q3.Cells.Item($r, 1).Value = $i

  $q3.Cells.Item($r, 2).NumberFormat = "@"
  $q3.Cells.Item($r, 2).Value = $fund[0]

  $q3.Cells.Item($r, 3).NumberFormat = "@"
  $q3.Cells.Item($r, 3).Value = $fund[1]

  $q3.Cells.Item($r, 4).NumberFormat = "@"
  $q3.Cells.Item($r, 4).Value = $fund[2]

  $q3.Cells.Item($r, 5).NumberFormat = "@"
  $q3.Cells.Item($r, 5).Value = $fund[3]

  $q3.Cells.Item($r, 6).NumberFormat = "@"
  $q3.Cells.Item($r, 6).Value = $fund[4]

  $q3.Cells.Item($r, 7).NumberFormat = "@"
  $q3.Cells.Item($r, 7).Value = $fund[5]

  $q3.Cells.Item($r, 8).Value = $fund[6]
}

$q3.Range("A1").Select()
$summary.Select()
